# Update sheet with newly computed TPM-derived values (K-T columns, rows 2-5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.110264333333333
$ws.Range("N2").Value = 6.330793
$ws.Range("O2").Value = 0.3832041185227171
$ws.Range("P2").Value = 0.3832041185227171
$ws.Range("Q2").Value = 0.1149397674436667
$ws.Range("R2").Value = 1.034457906993
$ws.Range("S2").Value = 0.3832041185227171
$ws.Range("T2").Value = 0.3832041185227171

# Row 3
$ws.Range("O3").Value = 0.0946183755984393
$ws.Range("P3").Value = 0.0946183755984393
$ws.Range("S3").Value = 0.0946183755984393
$ws.Range("T3").Value = 0.0946183755984393

# Row 4
$ws.Range("M4").Value = 2.065388333333333
$ws.Range("N4").Value = 6.196165
$ws.Range("O4").Value = 0.3750550597762889
$ws.Range("P4").Value = 0.3750550597762889
$ws.Range("Q4").Value = 0.1124955063516667
$ws.Range("R4").Value = 1.012459557165
$ws.Range("S4").Value = 0.3750550597762889
$ws.Range("T4").Value = 0.3750550597762889

# Row 5
$ws.Range("M5").Value = 0.8101876666666666
$ws.Range("N5").Value = 2.430563
$ws.Range("O5").Value = 0.1471224461025547
$ws.Range("P5").Value = 0.1471224461025547
$ws.Range("Q5").Value = 0.04412849164033333
$ws.Range("R5").Value = 0.3971564247629999
$ws.Range("S5").Value = 0.1471224461025547
$ws.Range("T5").Value = 0.1471224461025547
